$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the old "Température" column (A) to the new column F, -----------
# --- keeping the original per-row formulas (=22.4, =30, =35, ...) ---------
$ws.Range("F2").Formula = '=22.4'
$ws.Range("F3").Formula = '=30'
$ws.Range("F4").Formula = '=35'
$ws.Range("F5").Formula = '=40'
$ws.Range("F6").Formula = '=45'
$ws.Range("F7").Formula = '=50'
$ws.Range("F8").Formula = '=55'
$ws.Range("F9").Formula = '=60'
$ws.Range("F10").Formula = '=65'
$ws.Range("F11").Formula = '=70'
$ws.Range("F12").Formula = '=75'
$ws.Range("F13").Formula = '=80'
$ws.Range("F14").Formula = '=85'

# --- new O8 helper cell (was E5) -------------------------------------------
$ws.Range("O8").Formula = '=10^-4'

# --- column A becomes "DeltaTemp" = F{r} - $F$2 -----------------------------
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("A$r").Formula = "=F$r-`$F`$2"
}

# --- column B now multiplies by $O$8 instead of $E$5 ------------------------
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("B$r").Formula = "=C$r*`$O`$8"
}

# --- C2 resistance value changed from 94 to 95 ------------------------------
$ws.Range("C2").Formula = '=95'

# --- column E = "Resistance-" = B{r}/$B$2 - 1 -------------------------------
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("E$r").Formula = "=B$r/`$B`$2-1"
}

# --- column G = "AllDeltaTemp" ----------------------------------------------
$ws.Range("G2").Formula = '=-22.4'
for ($r = 3; $r -le 14; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=G$prev+10"
}

# --- column H = "AllResistivite" --------------------------------------------
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("H$r").Formula = "=`$D`$2*(1+0.0029*G$r)"
}

# --- header row --------------------------------------------------------------
$ws.Range("A1").Value = "DeltaTemp"
$ws.Range("E1").Value = "Resistance-"
$ws.Range("F1").Value = "Température"
$ws.Range("G1").Value = "AllDeltaTemp"
$ws.Range("H1").Value = "AllResistivite"

# --- cosmetics: column widths (character units -> internal width ~ +5/6) ----
$ws.Columns("F").ColumnWidth = 11.166666666666666
$ws.Columns("G").ColumnWidth = 11.5
$ws.Columns("H").ColumnWidth = 11.666666666666666

# --- cosmetics: view / selection ---------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 104
$ws.Range("G15:H17").Select()
